$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 119. This pushes the existing rows 119-188
# down to 120-189 (preserving all of their data), matching the diff which
# shows row N's old content becoming row N+1's content, with a brand new
# row 119 added at the top of the shifted block.
$ws.Rows("119:119").Insert()

# Populate the newly inserted row 119 with the new data record.
$ws.Range("A119").Value = 4
$ws.Range("B119").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C119").Value = "Los Lagos"
$ws.Range("D119").Value = 44572
$ws.Range("E119").Value = 10
$ws.Range("F119").Value = "Fruta"
$ws.Range("G119").Value = 100108
$ws.Range("H119").Value = "Tropicales y subtropicales"
$ws.Range("I119").Value = 100108005
$ws.Range("J119").Value = "Piña"
$ws.Range("K119").Value = "Caramelo"
$ws.Range("L119").Value = "Tercera"
$ws.Range("M119").Value = 200
$ws.Range("N119").Value = 19000
$ws.Range("O119").Value = 20000
$ws.Range("P119").Value = 19500
$ws.Range("Q119").Value = "$/caja 16 unidades"
$ws.Range("R119").Value = "Ecuador"
$ws.Range("S119").Value = 1219
$ws.Range("T119").Value = 16
